$wb = $excel.ActiveWorkbook

# --- ALC row 17 (hunk 0) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3626.1667
$ws.Range("J17").Value = 3721.8235
$ws.Range("L17").Value = 11165.4705
$ws.Range("N17").Value = -11501.4705

# --- ALC row 40 (hunk 1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2163.32
$ws.Range("J40").Value = 2170.8572
$ws.Range("L40").Value = 2170.8572
$ws.Range("N40").Value = -2520.8572

# --- ALC row 111 (hunk 2) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1060.7142
$ws.Range("I111").Value = 442.33334
$ws.Range("J111").Value = 1524.5
$ws.Range("K111").Value = 1327.00002
$ws.Range("L111").Value = 4573.5
$ws.Range("M111").Value = 1739.99998
$ws.Range("N111").Value = -10707.5

# --- ALC row 137 (hunk 3) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1419.6666
$ws.Range("I137").Value = 1419.6666
$ws.Range("K137").Value = 4258.9998
$ws.Range("M137").Value = -1708.9998

# --- ALC row 138 (hunk 4) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2619.9143
$ws.Range("J138").Value = 2215.1738
$ws.Range("L138").Value = 6645.5214
$ws.Range("N138").Value = -16925.5214

# --- ARM row 11 (hunk 5) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 1002400
$ws.Range("I11").Value = 1002400
$ws.Range("K11").Value = 1002400
$ws.Range("M11").Value = -1002256

# --- ARM row 32 (hunk 6) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6912.7334
$ws.Range("I32").Value = 7049.357
$ws.Range("K32").Value = 7049.357
$ws.Range("M32").Value = -6762.357

# --- ARM row 45 (hunk 7) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5077.0557
$ws.Range("I45").Value = 5559.1333
$ws.Range("J45").Value = 2666.6667
$ws.Range("K45").Value = 5559.1333
$ws.Range("L45").Value = 2666.6667
$ws.Range("M45").Value = -5182.1333
$ws.Range("N45").Value = -3420.6667

# --- ARM row 61 (hunk 8) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6624.125
$ws.Range("I61").Value = 2499
$ws.Range("K61").Value = 2499
$ws.Range("M61").Value = -2287

# --- ARM row 101 (hunk 9) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 349000
$ws.Range("J101").Value = 349000
$ws.Range("L101").Value = 349000
$ws.Range("N101").Value = -355490

# --- ARM row 102 (hunk 10) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 969
$ws.Range("I102").Value = 969
$ws.Range("K102").Value = 969
$ws.Range("M102").Value = 653

# --- ARM row 136 (hunk 11) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6624.125
$ws.Range("I136").Value = 2499
$ws.Range("K136").Value = 7497
$ws.Range("M136").Value = -4947

# --- BSM row 55 (hunk 12) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 78000
$ws.Range("J55").Value = 78000
$ws.Range("L55").Value = 78000
$ws.Range("N55").Value = -78546

# --- BSM row 99 (hunk 13) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1320.8889
$ws.Range("I99").Value = 1320.8889
$ws.Range("K99").Value = 1320.8889
$ws.Range("M99").Value = 177.1111000000001

# --- BSM row 105 (hunk 14) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3463.1428
$ws.Range("I105").Value = 3463.1428
$ws.Range("K105").Value = 3463.1428
$ws.Range("M105").Value = -1716.1428

# --- BSM row 134 (hunk 15) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 996.5
$ws.Range("I134").Value = 996.5
$ws.Range("K134").Value = 2989.5
$ws.Range("M134").Value = -454.5

# --- CRP row 31 (hunk 16) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1623.25
$ws.Range("I31").Value = 1325.1666
$ws.Range("J31").Value = 1921.3334
$ws.Range("K31").Value = 1325.1666
$ws.Range("L31").Value = 1921.3334
$ws.Range("M31").Value = -1030.1666
$ws.Range("N31").Value = -2511.3334

# --- CRP row 34 (hunk 17) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1623.25
$ws.Range("I34").Value = 1325.1666
$ws.Range("J34").Value = 1921.3334
$ws.Range("K34").Value = 1325.1666
$ws.Range("L34").Value = 1921.3334
$ws.Range("M34").Value = -1123.1666
$ws.Range("N34").Value = -2325.3334

# --- CRP row 43 (hunk 18) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 13000
$ws.Range("J43").Value = 15000
$ws.Range("L43").Value = 15000
$ws.Range("N43").Value = -15368

# --- CRP row 101 (hunk 19) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H101").Value = 13000
$ws.Range("J101").Value = 15000
$ws.Range("L101").Value = 15000
$ws.Range("N101").Value = -21490

# --- CUL row 37 (hunk 20) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 63993.332
$ws.Range("J37").Value = 63993.332
$ws.Range("L37").Value = 191979.996
$ws.Range("N37").Value = -192203.996

# --- CUL row 129 (hunk 21) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1002780.2
$ws.Range("I129").Value = 1600
$ws.Range("K129").Value = 4800
$ws.Range("M129").Value = 200

# --- GSM row 122 (hunk 22) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3768.625
$ws.Range("I122").Value = 3809.8
$ws.Range("J122").Value = 3700
$ws.Range("K122").Value = 11429.4
$ws.Range("L122").Value = 11100
$ws.Range("M122").Value = -8979.400000000001
$ws.Range("N122").Value = -16000

# --- LTW row 16 (hunk 23) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2919.8
$ws.Range("I16").Value = 1866.3334
$ws.Range("K16").Value = 1866.3334
$ws.Range("M16").Value = -1696.3334

# --- LTW row 46 (hunk 24) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2881.1667
$ws.Range("I46").Value = 2857.4
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 2857.4
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -2669.4
$ws.Range("N46").Value = -3376

# --- LTW row 68 (hunk 25) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2433
$ws.Range("I68").Value = 1999.3334
$ws.Range("J68").Value = 2866.6667
$ws.Range("K68").Value = 1999.3334
$ws.Range("L68").Value = 2866.6667
$ws.Range("M68").Value = -1250.3334
$ws.Range("N68").Value = -4364.6667

# --- LTW row 71 (hunk 26) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2433
$ws.Range("I71").Value = 1999.3334
$ws.Range("J71").Value = 2866.6667
$ws.Range("K71").Value = 9996.666999999999
$ws.Range("L71").Value = 14333.3335
$ws.Range("M71").Value = -6252.666999999999
$ws.Range("N71").Value = -21821.3335

# --- LTW row 100 (hunk 27) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3814.0715
$ws.Range("I100").Value = 2988.7778
$ws.Range("K100").Value = 2988.7778
$ws.Range("M100").Value = -2447.7778

# --- LTW row 101 (hunk 28) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

# --- LTW row 122 (hunk 29) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7691.3
$ws.Range("I122").Value = 7965.4375
$ws.Range("K122").Value = 23896.3125
$ws.Range("M122").Value = -21446.3125

# --- LTW row 132 (hunk 30) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2543.8572
$ws.Range("I132").Value = 2055.3635
$ws.Range("K132").Value = 6166.0905
$ws.Range("M132").Value = -3636.0905

# --- WVR row 107 (hunk 31) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1990.5
$ws.Range("J107").Value = 1980
$ws.Range("L107").Value = 5940
$ws.Range("N107").Value = -9780

# --- WVR row 122 (hunk 32) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 924
$ws.Range("I122").Value = 937.25
$ws.Range("K122").Value = 2811.75
$ws.Range("M122").Value = -361.75

# --- WVR row 136 (hunk 33) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3634.5417
$ws.Range("I136").Value = 3495
$ws.Range("J136").Value = 3973.4285
$ws.Range("K136").Value = 10485
$ws.Range("L136").Value = 11920.2855
$ws.Range("M136").Value = -7935
$ws.Range("N136").Value = -17020.2855
